# Fixed y polarization and added values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Flip sign of "3d y" (col B) and col D values for rows 2-8 ---
$ws.Range("B2").Value = -0.35625000000000001
$ws.Range("D2").Value = 0.74099999999999999

$ws.Range("B3").Value = -1.0069999999999999
$ws.Range("D3").Value = 1.0545

$ws.Range("B4").Value = -0.06175
$ws.Range("D4").Value = 0.30399999999999999

$ws.Range("B5").Value = 0.07125
$ws.Range("D5").Value = 0.076

$ws.Range("B6").Value = 0.69350000000000001
$ws.Range("D6").Value = -1.0069999999999999

$ws.Range("B7").Value = 1.4724999999999999
$ws.Range("D7").Value = -0.96250000000000002

$ws.Range("B8").Value = 1.2017500000000001
$ws.Range("D8").Value = -0.78374999999999995

# --- Remove the now-stale "3d est" columns J:L for rows 2-8 ---
$ws.Range("J2:L8").ClearContents()

# --- Add two new data rows (9 and 10) ---
$ws.Range("A9").Value = -0.371
$ws.Range("B9").Value = 0.25175
$ws.Range("C9").Value = 0.5285
$ws.Range("D9").Value = -0.01425
$ws.Range("F9").Value = 15000
$ws.Range("G9").Value = 5000
$ws.Range("H9").Value = 4880

$ws.Range("A10").Value = -1.305
$ws.Range("B10").Value = 1.33475
$ws.Range("C10").Value = 0.8715
$ws.Range("D10").Value = -0.20425
$ws.Range("F10").Value = 12500
$ws.Range("G10").Value = 7500
$ws.Range("H10").Value = 4880

# --- Update the sheet's selection to match the committed state ---
$ws.Range("J2:L10").Select()
